$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset formatting on the amount/quantity/instrument cells in the
#     existing rows so they fall back to the default (General) style,
#     matching the cleanup that was done alongside the new test row.
$ws.Range("D2:E3").Style = "Normal"
$ws.Range("H2:H3").Style = "Normal"

# --- Build the new row 4 by basing each cell on the matching cell from
#     row 2/3 (so it inherits the same look/number format) and then
#     overwriting the value/text.
$ws.Range("B2").Copy($ws.Range("B4"))
$ws.Range("C2").Copy($ws.Range("C4"))
$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("I2").Copy($ws.Range("I4"))

$ws.Range("A4").Value = "Fund X"
$ws.Range("B4").Value = "Ego Pvt Ltd"
$ws.Range("C4").Value = 44880
$ws.Range("D4").Value = 300000
$ws.Range("E4").Value = 150000
$ws.Range("F4").Value = "Test"
$ws.Range("H4").Value = "Equity"
$ws.Range("I4").Value = "INR"

$ws.Range("E4").Select()
